$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '26.950.30'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.07%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.675.31'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +1.22%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '215.17'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.04%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.516'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.77%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +0.21%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.0620'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.39%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '20.21'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.56%  '
$ws.Range('E11').Value = '  +1.02%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.909.75'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +1.08%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.663.30'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.79%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.09'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.03%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.526'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.39%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '65.67'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.51%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '26.951.25'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.10%  '
$ws.Range('B18').Value = 'BitcoinCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '235.13'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.67%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '8.08'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +4.09%  '
$ws.Range('E20').Value = '  -0.72%  '
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.45'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.32%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.17'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.25%  '
$ws.Range('E24').Value = '  -1.73%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '145.58'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.26%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.17'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.64%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '16.04'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.07%  '
$ws.Range('E28').Value = '  -1.36%  '
$ws.Range('E29').Value = '  +0.18%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0498'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +0.07%  '
$ws.Range('E31').Value = '  +0.16%  '
$ws.Range('E32').Value = '  +1.00%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.478.99'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -4.90%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.15'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +2.19%  '
$ws.Range('E35').Value = '  +2.56%  '
$ws.Range('E36').Value = '  +0.10%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.583'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.22%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.897'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.00%  '
$ws.Range('E39').Value = '  +0.80%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.05'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +7.12%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.83'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -3.42%  '
$ws.Range('E42').Value = '  +0.05%  '
$ws.Range('E43').Value = '  +2.40%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '66.87'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.78%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.816.02'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.03%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.780'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.83%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '90.53'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.37%  '
$ws.Range('E48').Value = '  +0.21%  '
$ws.Range('E49').Value = '  +1.73%  '
$ws.Range('E50').Value = '  +0.34%  '
$ws.Range('E51').Value = '  +0.16%  '
